$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 81.03634266666667
$ws.Range("H2").Value = 243.109028
$ws.Range("I2").Value = 0.1632931649012984
$ws.Range("J2").Value = 0.1632931649012984
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 1915.358598916
$ws.Range("R2").Value = 17238.22739024401
$ws.Range("S2").Value = 0.01114940352166862
$ws.Range("T2").Value = 0.01114940352166862

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 81.03634266666667
$ws.Range("H3").Value = 243.109028
$ws.Range("I3").Value = 0.1632931649012984
$ws.Range("J3").Value = 0.1632931649012984
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 14690.94820157677
$ws.Range("R3").Value = 132218.5338141909
$ws.Range("S3").Value = 0.08551678505947212
$ws.Range("T3").Value = 0.0855167850594721

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 81.03634266666667
$ws.Range("H4").Value = 243.109028
$ws.Range("I4").Value = 0.1632931649012984
$ws.Range("J4").Value = 0.1632931649012984
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 9003.181943121877
$ws.Range("R4").Value = 81028.6374880969
$ws.Range("S4").Value = 0.05240799739520142
$ws.Range("T4").Value = 0.05240799739520142

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 81.03634266666667
$ws.Range("H5").Value = 243.109028
$ws.Range("I5").Value = 0.1632931649012984
$ws.Range("J5").Value = 0.1632931649012984
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 2442.681664430811
$ws.Range("R5").Value = 21984.1349798773
$ws.Range("S5").Value = 0.01421897892495621
$ws.Range("T5").Value = 0.01421897892495621

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 178.365814
$ws.Range("H6").Value = 535.097442
$ws.Range("I6").Value = 0.3594179761796791
$ws.Range("J6").Value = 0.3594179761796791
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 4215.8182903543
$ws.Range("R6").Value = 37942.36461318871
$ws.Range("S6").Value = 0.02454050083352178
$ws.Range("T6").Value = 0.02454050083352178

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 178.365814
$ws.Range("H7").Value = 535.097442
$ws.Range("I7").Value = 0.3594179761796791
$ws.Range("J7").Value = 0.3594179761796791
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 32335.65148892055
$ws.Range("R7").Value = 291020.8634002849
$ws.Range("S7").Value = 0.1882275344105582
$ws.Range("T7").Value = 0.1882275344105581

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 178.365814
$ws.Range("H8").Value = 535.097442
$ws.Range("I8").Value = 0.3594179761796791
$ws.Range("J8").Value = 0.3594179761796791
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 19816.53938258972
$ws.Range("R8").Value = 178348.8544433074
$ws.Range("S8").Value = 0.1153531219190878
$ws.Range("T8").Value = 0.1153531219190878

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 178.365814
$ws.Range("H9").Value = 535.097442
$ws.Range("I9").Value = 0.3594179761796791
$ws.Range("J9").Value = 0.3594179761796791
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 5376.487747124017
$ws.Range("R9").Value = 48388.38972411615
$ws.Range("S9").Value = 0.03129681901651129
$ws.Range("T9").Value = 0.03129681901651129

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 146.2303303333333
$ws.Range("H10").Value = 438.6909910000001
$ws.Range("I10").Value = 0.2946630198121519
$ws.Range("J10").Value = 0.2946630198121519
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 3456.270500488496
$ws.Range("R10").Value = 31106.43450439647
$ws.Range("S10").Value = 0.02011913304996513
$ws.Range("T10").Value = 0.02011913304996513

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 146.2303303333333
$ws.Range("H11").Value = 438.6909910000001
$ws.Range("I11").Value = 0.2946630198121519
$ws.Range("J11").Value = 0.2946630198121519
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 26509.86135027194
$ws.Range("R11").Value = 238588.7521524475
$ws.Range("S11").Value = 0.1543153024529958
$ws.Range("T11").Value = 0.1543153024529958

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 146.2303303333333
$ws.Range("H12").Value = 438.6909910000001
$ws.Range("I12").Value = 0.2946630198121519
$ws.Range("J12").Value = 0.2946630198121519
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 16246.26959053714
$ws.Range("R12").Value = 146216.4263148342
$ws.Range("S12").Value = 0.09457039297457241
$ws.Range("T12").Value = 0.09457039297457241

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 146.2303303333333
$ws.Range("H13").Value = 438.6909910000001
$ws.Range("I13").Value = 0.2946630198121519
$ws.Range("J13").Value = 0.2946630198121519
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 4407.826598963993
$ws.Range("R13").Value = 39670.43939067594
$ws.Range("S13").Value = 0.02565819133461861
$ws.Range("T13").Value = 0.02565819133461861

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 90.630432
$ws.Range("H14").Value = 271.891296
$ws.Range("I14").Value = 0.1826258391068707
$ws.Range("J14").Value = 0.1826258391068707
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 2142.122553194592
$ws.Range("R14").Value = 19279.10297875133
$ws.Range("S14").Value = 0.01246940847105622
$ws.Range("T14").Value = 0.01246940847105622

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 90.630432
$ws.Range("H15").Value = 271.891296
$ws.Range("I15").Value = 0.1826258391068707
$ws.Range("J15").Value = 0.1826258391068707
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 16430.24522312506
$ws.Range("R15").Value = 147872.2070081255
$ws.Range("S15").Value = 0.09564132484447806
$ws.Range("T15").Value = 0.09564132484447804

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 90.630432
$ws.Range("H16").Value = 271.891296
$ws.Range("I16").Value = 0.1826258391068707
$ws.Range("J16").Value = 0.1826258391068707
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 10069.09050962602
$ws.Range("R16").Value = 90621.81458663414
$ws.Range("S16").Value = 0.05861270743325064
$ws.Range("T16").Value = 0.05861270743325064

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 90.630432
$ws.Range("H17").Value = 271.891296
$ws.Range("I17").Value = 0.1826258391068707
$ws.Range("J17").Value = 0.1826258391068707
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 2731.876676572992
$ws.Range("R17").Value = 24586.89008915692
$ws.Range("S17").Value = 0.01590239835808578
$ws.Range("T17").Value = 0.01590239835808578

